$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "source_data"

# Update existing values / add new values per diff
# Row 2 (Coal): co2_rate (K2) and land_rate (M2)
$ws.Range("K2").Value = 96
$ws.Range("M2").Value = 8000

# Row 3 (Natural Gas): co2_rate (K3) and land_rate (M3)
$ws.Range("K3").Value = 57
$ws.Range("M3").Value = 2000

# Row 4 (Advanced Nuclear): fuel_cost (G4) changes, land_rate (M4) added
$ws.Range("G4").Value = 1
$ws.Range("M4").Value = 1000

# Row 5 (Onshore Wind): land_rate (M5) added
$ws.Range("M5").Value = 5

# Row 6 (new: Rooftop Solar PV): name (A6) and land_rate (M6)
$ws.Range("A6").Value = "Rooftop Solar PV"
$ws.Range("M6").Value = 10

# Update selection to match diff
$ws.Range("H12").Select()
